# Sheet4 ("存款" / deposit) gains new columns: bank/deposit_type/currency stay
# in B:D, but the amount moves into F (matching the other sheets' "total"
# column), and new metadata columns (category/normal/date/legislator_*/index)
# are appended in G:M, matching the shared layout used by sheet1-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- First, extend formatting to the new cells (H1:M1 use the bold/boxed
# header style already on B1; G2:M12 use the plain boxed data style already
# on B2). F1:F12 and G1 already existed with the correct style, so only
# H1:M1 and G2:M12 need their format copied in from a neighbour. ---

$ws.Range("E1").Copy()
$ws.Range("H1:M1").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("G2:M12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Header row (row 1) ---
$ws.Cells.Item(1,2).Value = "bank"
$ws.Cells.Item(1,3).Value = "deposit_type"
$ws.Cells.Item(1,4).Value = "currency"
$ws.Cells.Item(1,5).Value = "owner"
$ws.Cells.Item(1,6).Value = "total"
$ws.Cells.Item(1,7).Value = "property_category"
$ws.Cells.Item(1,8).Value = "category"
$ws.Cells.Item(1,9).Value = "date"
$ws.Cells.Item(1,10).Value = "legislator_name"
$ws.Cells.Item(1,11).Value = "legislator_id"
$ws.Cells.Item(1,12).Value = "source_file"
$ws.Cells.Item(1,13).Value = "index"

# --- Data rows (row 2-12) ---
$ws.Cells.Item(2,1).Value = 80
$ws.Cells.Item(2,2).Value = "臺灣銀行群賢分行"
$ws.Cells.Item(2,3).Value = "活期儲蓄存款"
$ws.Cells.Item(2,4).Value = "新臺幣"
$ws.Cells.Item(2,5).Value = "管碧玲"
$ws.Cells.Item(2,6).Value = 3273822
$ws.Cells.Item(2,7).Value = "deposit"
$ws.Cells.Item(2,8).Value = "normal"
$ws.Cells.Item(2,9).Value = "2013-12-30"
$ws.Cells.Item(2,10).Value = "管碧玲"
$ws.Cells.Item(2,11).Value = 1374
$ws.Cells.Item(2,12).Value = "tmpb8981"
$ws.Cells.Item(2,13).Value = 80

$ws.Cells.Item(3,1).Value = 81
$ws.Cells.Item(3,2).Value = "合作金庫商業銀行西門分行"
$ws.Cells.Item(3,3).Value = "活期儲蓄存款"
$ws.Cells.Item(3,4).Value = "新臺幣"
$ws.Cells.Item(3,5).Value = "管碧玲"
$ws.Cells.Item(3,6).Value = 5384
$ws.Cells.Item(3,7).Value = "deposit"
$ws.Cells.Item(3,8).Value = "normal"
$ws.Cells.Item(3,9).Value = "2013-12-30"
$ws.Cells.Item(3,10).Value = "管碧玲"
$ws.Cells.Item(3,11).Value = 1374
$ws.Cells.Item(3,12).Value = "tmpb8981"
$ws.Cells.Item(3,13).Value = 81

$ws.Cells.Item(4,1).Value = 82
$ws.Cells.Item(4,2).Value = "高雄銀行市府分行"
$ws.Cells.Item(4,3).Value = "活期儲蓄存款"
$ws.Cells.Item(4,4).Value = "新臺幣"
$ws.Cells.Item(4,5).Value = "管碧玲"
$ws.Cells.Item(4,6).Value = 1667936
$ws.Cells.Item(4,7).Value = "deposit"
$ws.Cells.Item(4,8).Value = "normal"
$ws.Cells.Item(4,9).Value = "2013-12-30"
$ws.Cells.Item(4,10).Value = "管碧玲"
$ws.Cells.Item(4,11).Value = 1374
$ws.Cells.Item(4,12).Value = "tmpb8981"
$ws.Cells.Item(4,13).Value = 82

$ws.Cells.Item(5,1).Value = 83
$ws.Cells.Item(5,2).Value = "國泰世華商業銀行中正分行"
$ws.Cells.Item(5,3).Value = "活期儲蓄存款"
$ws.Cells.Item(5,4).Value = "新臺幣"
$ws.Cells.Item(5,5).Value = "管碧玲"
$ws.Cells.Item(5,6).Value = 433548
$ws.Cells.Item(5,7).Value = "deposit"
$ws.Cells.Item(5,8).Value = "normal"
$ws.Cells.Item(5,9).Value = "2013-12-30"
$ws.Cells.Item(5,10).Value = "管碧玲"
$ws.Cells.Item(5,11).Value = 1374
$ws.Cells.Item(5,12).Value = "tmpb8981"
$ws.Cells.Item(5,13).Value = 83

$ws.Cells.Item(6,1).Value = 84
$ws.Cells.Item(6,2).Value = "高雄市府郵局(第19支局)"
$ws.Cells.Item(6,3).Value = "中華郵政存簿儲金"
$ws.Cells.Item(6,4).Value = "新臺幣"
$ws.Cells.Item(6,5).Value = "管碧玲"
$ws.Cells.Item(6,6).Value = 493008
$ws.Cells.Item(6,7).Value = "deposit"
$ws.Cells.Item(6,8).Value = "normal"
$ws.Cells.Item(6,9).Value = "2013-12-30"
$ws.Cells.Item(6,10).Value = "管碧玲"
$ws.Cells.Item(6,11).Value = 1374
$ws.Cells.Item(6,12).Value = "tmpb8981"
$ws.Cells.Item(6,13).Value = 84

$ws.Cells.Item(7,1).Value = 86
$ws.Cells.Item(7,2).Value = "合作金庫商業銀行十全分行"
$ws.Cells.Item(7,3).Value = "活期存款"
$ws.Cells.Item(7,4).Value = "新臺幣"
$ws.Cells.Item(7,5).Value = "管碧玲"
$ws.Cells.Item(7,6).Value = 20310
$ws.Cells.Item(7,7).Value = "deposit"
$ws.Cells.Item(7,8).Value = "normal"
$ws.Cells.Item(7,9).Value = "2013-12-30"
$ws.Cells.Item(7,10).Value = "管碧玲"
$ws.Cells.Item(7,11).Value = 1374
$ws.Cells.Item(7,12).Value = "tmpb8981"
$ws.Cells.Item(7,13).Value = 86

$ws.Cells.Item(8,1).Value = 87
$ws.Cells.Item(8,2).Value = "台北富邦商業銀行北投分行"
$ws.Cells.Item(8,3).Value = "活期儲蓄存款"
$ws.Cells.Item(8,4).Value = "新臺幣"
$ws.Cells.Item(8,5).Value = "管碧玲"
$ws.Cells.Item(8,6).Value = 1238
$ws.Cells.Item(8,7).Value = "deposit"
$ws.Cells.Item(8,8).Value = "normal"
$ws.Cells.Item(8,9).Value = "2013-12-30"
$ws.Cells.Item(8,10).Value = "管碧玲"
$ws.Cells.Item(8,11).Value = 1374
$ws.Cells.Item(8,12).Value = "tmpb8981"
$ws.Cells.Item(8,13).Value = 87

$ws.Cells.Item(9,1).Value = 88
$ws.Cells.Item(9,2).Value = "陽信商業銀行大屯分行"
$ws.Cells.Item(9,3).Value = "活期儲蓄存款"
$ws.Cells.Item(9,4).Value = "新臺幣"
$ws.Cells.Item(9,5).Value = "管碧玲"
$ws.Cells.Item(9,6).Value = 9792
$ws.Cells.Item(9,7).Value = "deposit"
$ws.Cells.Item(9,8).Value = "normal"
$ws.Cells.Item(9,9).Value = "2013-12-30"
$ws.Cells.Item(9,10).Value = "管碧玲"
$ws.Cells.Item(9,11).Value = 1374
$ws.Cells.Item(9,12).Value = "tmpb8981"
$ws.Cells.Item(9,13).Value = 88

$ws.Cells.Item(10,1).Value = 89
$ws.Cells.Item(10,2).Value = "合作金庫商業銀行長春分行"
$ws.Cells.Item(10,3).Value = "活期存款"
$ws.Cells.Item(10,4).Value = "新臺幣"
$ws.Cells.Item(10,5).Value = "管碧玲"
$ws.Cells.Item(10,6).Value = 31205
$ws.Cells.Item(10,7).Value = "deposit"
$ws.Cells.Item(10,8).Value = "normal"
$ws.Cells.Item(10,9).Value = "2013-12-30"
$ws.Cells.Item(10,10).Value = "管碧玲"
$ws.Cells.Item(10,11).Value = 1374
$ws.Cells.Item(10,12).Value = "tmpb8981"
$ws.Cells.Item(10,13).Value = 89

$ws.Cells.Item(11,1).Value = 90
$ws.Cells.Item(11,2).Value = "臺灣銀行北投分行"
$ws.Cells.Item(11,3).Value = "活期儲蓄存款"
$ws.Cells.Item(11,4).Value = "新臺幣"
$ws.Cells.Item(11,5).Value = "管碧玲"
$ws.Cells.Item(11,6).Value = 108527
$ws.Cells.Item(11,7).Value = "deposit"
$ws.Cells.Item(11,8).Value = "normal"
$ws.Cells.Item(11,9).Value = "2013-12-30"
$ws.Cells.Item(11,10).Value = "管碧玲"
$ws.Cells.Item(11,11).Value = 1374
$ws.Cells.Item(11,12).Value = "tmpb8981"
$ws.Cells.Item(11,13).Value = 90

$ws.Cells.Item(12,1).Value = 91
$ws.Cells.Item(12,2).Value = "合作金庫商業銀行營業部"
$ws.Cells.Item(12,3).Value = "活期存款"
$ws.Cells.Item(12,4).Value = "美金"
$ws.Cells.Item(12,5).Value = "管碧玲"
$ws.Cells.Item(12,6).Value = 118465.95
$ws.Cells.Item(12,7).Value = "deposit"
$ws.Cells.Item(12,8).Value = "normal"
$ws.Cells.Item(12,9).Value = "2013-12-30"
$ws.Cells.Item(12,10).Value = "管碧玲"
$ws.Cells.Item(12,11).Value = 1374
$ws.Cells.Item(12,12).Value = "tmpb8981"
$ws.Cells.Item(12,13).Value = 91
